# Refining timings & blueprints
# Update the "Lao" (D) and "Total" (G) columns for the "Meteors" (row 2)
# and "Successes" (row 3) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2290
$ws.Range("G2").Value = 2290

$ws.Range("D3").Value = 22
$ws.Range("G3").Value = 22
